$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Prepend a "<ab>" run (Courier New, 7f6000, 9pt) right before the
#    "Pour un " run that starts the paragraph currently missing its opening
#    <ab> tag.
# ---------------------------------------------------------------------------
$rTarget1 = $d.Content
$find1 = $rTarget1.Find
$find1.ClearFormatting()
[void]$find1.Execute("Pour un ", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$insPos1 = $rTarget1.Start

# Create a small placeholder run at the insertion point so we have a stable,
# correctly-ordered target range to overwrite with formatted text.
$rPlaceholder1 = $d.Range($insPos1, $insPos1)
$rPlaceholder1.InsertBefore("ZZZZ")

# Borrow the fully-specified run formatting from an existing "<ab>" run
# elsewhere in the document (keeps rFonts/color/sz/szCs/rtl identical).
$rSrc1 = $d.Content
$findSrc1 = $rSrc1.Find
$findSrc1.ClearFormatting()
[void]$findSrc1.Execute("<ab>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$ft1 = $rSrc1.FormattedText

$rTargetFinal1 = $d.Range($insPos1, $insPos1 + 4)
$rTargetFinal1.FormattedText = $ft1

# ---------------------------------------------------------------------------
# 2) Append a "</ab>" run (Courier New, 7f6000, 9pt) right after the
#    "La chaulde<lb/>" content, before the paragraph's trailing empty run.
# ---------------------------------------------------------------------------
$rTarget2 = $d.Content
$find2 = $rTarget2.Find
$find2.ClearFormatting()
[void]$find2.Execute("La chaulde<lb/>", $false, $false, $false, $false, `
                      $false, $true, 1, $false, "", 0)
$insPos2 = $rTarget2.End

$rPlaceholder2 = $d.Range($insPos2, $insPos2)
$rPlaceholder2.InsertBefore("WWWWW")

# Borrow formatting from an existing "</ab>" run (the one that already
# terminates the "Pour un ..." paragraph).
$rSrc2 = $d.Content
$findSrc2 = $rSrc2.Find
$findSrc2.ClearFormatting()
[void]$findSrc2.Execute("</ab>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$ft2 = $rSrc2.FormattedText

$rTargetFinal2 = $d.Range($insPos2, $insPos2 + 5)
$rTargetFinal2.FormattedText = $ft2

# ---------------------------------------------------------------------------
# 3) Section page margins: add an explicit footer distance of 720 twips
#    (36pt) to the section's pgMar.
# ---------------------------------------------------------------------------
$ps = $d.PageSetup
$ps.FooterDistance = 36
